$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Role" column header and values to reflect the new
# PI/Sub I terminology (previously Principal/Sub).
$ws.Range("A1").Value = "Role (PI/Sub I)"
$ws.Range("A2").Value = "PI"
$ws.Range("A3").Value = "Sub I"
$ws.Range("A4").Value = "Sub I"
$ws.Range("A5").Value = "Sub I"

# Move the active selection to A6, as in the target workbook.
$ws.Range("A6").Select()
